# Auto-generated update of market price / profit columns (H-N) per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 1503.3636
$ws.Cells.Item(129, 9).Value = 1634.5
$ws.Cells.Item(129, 10).Value = 1461.4
$ws.Cells.Item(129, 11).Value = 4903.5
$ws.Cells.Item(129, 12).Value = 4384.200000000001
$ws.Cells.Item(129, 13).Value = 96.5
$ws.Cells.Item(129, 14).Value = -14384.2

$ws.Cells.Item(130, 8).Value = 43160
$ws.Cells.Item(130, 10).Value = 43160
$ws.Cells.Item(130, 12).Value = 43160
$ws.Cells.Item(130, 14).Value = -53200

$ws.Cells.Item(132, 8).Value = 31377.219
$ws.Cells.Item(132, 9).Value = 4367.643
$ws.Cells.Item(132, 10).Value = 220444.25
$ws.Cells.Item(132, 11).Value = 13102.929
$ws.Cells.Item(132, 12).Value = 661332.75
$ws.Cells.Item(132, 13).Value = -10572.929
$ws.Cells.Item(132, 14).Value = -666392.75

$ws.Cells.Item(138, 8).Value = 2117.679
$ws.Cells.Item(138, 9).Value = 1259.9796
$ws.Cells.Item(138, 10).Value = 3431.0312
$ws.Cells.Item(138, 11).Value = 3779.9388
$ws.Cells.Item(138, 12).Value = 10293.0936
$ws.Cells.Item(138, 13).Value = 1360.0612
$ws.Cells.Item(138, 14).Value = -20573.0936


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8165.2827
$ws.Cells.Item(32, 9).Value = 6870.988
$ws.Cells.Item(32, 11).Value = 6870.988
$ws.Cells.Item(32, 13).Value = -6583.988

$ws.Cells.Item(61, 8).Value = 3428.182
$ws.Cells.Item(61, 9).Value = 2300
$ws.Cells.Item(61, 10).Value = 4368.3335
$ws.Cells.Item(61, 11).Value = 2300
$ws.Cells.Item(61, 12).Value = 4368.3335
$ws.Cells.Item(61, 13).Value = -2088
$ws.Cells.Item(61, 14).Value = -4792.3335

$ws.Cells.Item(101, 8).Value = 48546
$ws.Cells.Item(101, 10).Value = 48546
$ws.Cells.Item(101, 12).Value = 48546
$ws.Cells.Item(101, 14).Value = -55036

$ws.Cells.Item(104, 8).Value = 21475.334
$ws.Cells.Item(104, 10).Value = 21475.334
$ws.Cells.Item(104, 12).Value = 21475.334
$ws.Cells.Item(104, 14).Value = -28463.334

$ws.Cells.Item(122, 8).Value = 1860.5883
$ws.Cells.Item(122, 9).Value = 1975.3846
$ws.Cells.Item(122, 11).Value = 5926.1538
$ws.Cells.Item(122, 13).Value = -3476.1538

$ws.Cells.Item(131, 8).Value = 45771.832
$ws.Cells.Item(131, 10).Value = 45771.832
$ws.Cells.Item(131, 12).Value = 45771.832
$ws.Cells.Item(131, 14).Value = -55851.832

$ws.Cells.Item(136, 8).Value = 3428.182
$ws.Cells.Item(136, 9).Value = 2300
$ws.Cells.Item(136, 10).Value = 4368.3335
$ws.Cells.Item(136, 11).Value = 6900
$ws.Cells.Item(136, 12).Value = 13105.0005
$ws.Cells.Item(136, 13).Value = -4350
$ws.Cells.Item(136, 14).Value = -18205.0005


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 10792
$ws.Cells.Item(26, 9).Value = 5530.2856
$ws.Cells.Item(26, 11).Value = 5530.2856
$ws.Cells.Item(26, 13).Value = -5238.2856

$ws.Cells.Item(59, 8).Value = 33054.285
$ws.Cells.Item(59, 9).Value = 20000
$ws.Cells.Item(59, 10).Value = 38276
$ws.Cells.Item(59, 11).Value = 20000
$ws.Cells.Item(59, 12).Value = 38276
$ws.Cells.Item(59, 14).Value = -39970
$ws.Cells.Item(59, 13).Value = -19153

$ws.Cells.Item(95, 8).Value = 43892
$ws.Cells.Item(95, 10).Value = 43892
$ws.Cells.Item(95, 12).Value = 43892
$ws.Cells.Item(95, 14).Value = -49384

$ws.Cells.Item(96, 8).Value = 12394
$ws.Cells.Item(96, 9).Value = 1745.6
$ws.Cells.Item(96, 10).Value = 20000
$ws.Cells.Item(96, 11).Value = 1745.6
$ws.Cells.Item(96, 12).Value = 20000
$ws.Cells.Item(96, 13).Value = 1000.4
$ws.Cells.Item(96, 14).Value = -25492

$ws.Cells.Item(100, 8).Value = 44786.668
$ws.Cells.Item(100, 10).Value = 44786.668
$ws.Cells.Item(100, 12).Value = 44786.668
$ws.Cells.Item(100, 14).Value = -46950.668

$ws.Cells.Item(124, 8).Value = 48974.668
$ws.Cells.Item(124, 10).Value = 48974.668
$ws.Cells.Item(124, 12).Value = 48974.668
$ws.Cells.Item(124, 14).Value = -58794.668

$ws.Cells.Item(130, 8).Value = 45081.8
$ws.Cells.Item(130, 10).Value = 45081.8
$ws.Cells.Item(130, 12).Value = 45081.8
$ws.Cells.Item(130, 14).Value = -55121.8

$ws.Cells.Item(134, 8).Value = 230157.38
$ws.Cells.Item(134, 9).Value = 3168.6667
$ws.Cells.Item(134, 10).Value = 254477.6
$ws.Cells.Item(134, 11).Value = 9506.000100000001
$ws.Cells.Item(134, 12).Value = 763432.8
$ws.Cells.Item(134, 13).Value = -6971.000100000001
$ws.Cells.Item(134, 14).Value = -768502.8


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2600.46
$ws.Cells.Item(31, 9).Value = 1056.7354
$ws.Cells.Item(31, 10).Value = 3395.7122
$ws.Cells.Item(31, 11).Value = 1056.7354
$ws.Cells.Item(31, 12).Value = 3395.7122
$ws.Cells.Item(31, 13).Value = -761.7354
$ws.Cells.Item(31, 14).Value = -3985.7122

$ws.Cells.Item(34, 8).Value = 2600.46
$ws.Cells.Item(34, 9).Value = 1056.7354
$ws.Cells.Item(34, 10).Value = 3395.7122
$ws.Cells.Item(34, 11).Value = 1056.7354
$ws.Cells.Item(34, 12).Value = 3395.7122
$ws.Cells.Item(34, 13).Value = -854.7354
$ws.Cells.Item(34, 14).Value = -3799.7122

$ws.Cells.Item(38, 8).Value = 2500
$ws.Cells.Item(38, 9).Value = 2000
$ws.Cells.Item(38, 10).Value = 3000
$ws.Cells.Item(38, 11).Value = 2000
$ws.Cells.Item(38, 12).Value = 3000
$ws.Cells.Item(38, 13).Value = -1623
$ws.Cells.Item(38, 14).Value = -3754

$ws.Cells.Item(46, 8).Value = 2500
$ws.Cells.Item(46, 9).Value = 2000
$ws.Cells.Item(46, 10).Value = 3000
$ws.Cells.Item(46, 11).Value = 2000
$ws.Cells.Item(46, 12).Value = 3000
$ws.Cells.Item(46, 13).Value = -1789
$ws.Cells.Item(46, 14).Value = -3422

$ws.Cells.Item(100, 8).Value = 46945.332
$ws.Cells.Item(100, 10).Value = 46945.332
$ws.Cells.Item(100, 12).Value = 46945.332
$ws.Cells.Item(100, 14).Value = -49109.332

$ws.Cells.Item(111, 8).Value = 49247
$ws.Cells.Item(111, 10).Value = 49247
$ws.Cells.Item(111, 12).Value = 49247
$ws.Cells.Item(111, 14).Value = -57427

$ws.Cells.Item(131, 8).Value = 38326
$ws.Cells.Item(131, 10).Value = 38326
$ws.Cells.Item(131, 12).Value = 38326
$ws.Cells.Item(131, 14).Value = -48406

$ws.Cells.Item(132, 8).Value = 28592.404
$ws.Cells.Item(132, 9).Value = 1442.1555
$ws.Cells.Item(132, 10).Value = 203129.72
$ws.Cells.Item(132, 11).Value = 4326.4665
$ws.Cells.Item(132, 12).Value = 609389.16
$ws.Cells.Item(132, 13).Value = -1796.4665
$ws.Cells.Item(132, 14).Value = -614449.16

$ws.Cells.Item(134, 8).Value = 425870.53
$ws.Cells.Item(134, 9).Value = 1151
$ws.Cells.Item(134, 10).Value = 1753119.1
$ws.Cells.Item(134, 11).Value = 3453
$ws.Cells.Item(134, 12).Value = 5259357.300000001
$ws.Cells.Item(134, 13).Value = -918
$ws.Cells.Item(134, 14).Value = -5264427.300000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 4201.788
$ws.Cells.Item(5, 9).Value = 17100.666
$ws.Cells.Item(5, 10).Value = 1335.3704
$ws.Cells.Item(5, 11).Value = 51301.99800000001
$ws.Cells.Item(5, 12).Value = 4006.1112
$ws.Cells.Item(5, 13).Value = -51189.99800000001
$ws.Cells.Item(5, 14).Value = -4230.1112

$ws.Cells.Item(113, 8).Value = 4370.1113
$ws.Cells.Item(113, 9).Value = 13281.375
$ws.Cells.Item(113, 10).Value = 618
$ws.Cells.Item(113, 11).Value = 39844.125
$ws.Cells.Item(113, 12).Value = 1854
$ws.Cells.Item(113, 13).Value = -37674.125
$ws.Cells.Item(113, 14).Value = -6194

$ws.Cells.Item(135, 8).Value = 4201.788
$ws.Cells.Item(135, 9).Value = 17100.666
$ws.Cells.Item(135, 10).Value = 1335.3704
$ws.Cells.Item(135, 11).Value = 153905.994
$ws.Cells.Item(135, 12).Value = 12018.3336
$ws.Cells.Item(135, 13).Value = -151370.994
$ws.Cells.Item(135, 14).Value = -17088.3336


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 78008.89999999999
$ws.Cells.Item(22, 9).Value = 10008
$ws.Cells.Item(22, 10).Value = 85564.55499999999
$ws.Cells.Item(22, 11).Value = 10008
$ws.Cells.Item(22, 12).Value = 85564.55499999999
$ws.Cells.Item(22, 13).Value = -9479
$ws.Cells.Item(22, 14).Value = -86622.55499999999

$ws.Cells.Item(104, 8).Value = 46462.75
$ws.Cells.Item(104, 10).Value = 46462.75
$ws.Cells.Item(104, 12).Value = 46462.75
$ws.Cells.Item(104, 14).Value = -53450.75

$ws.Cells.Item(113, 8).Value = 1591.45
$ws.Cells.Item(113, 9).Value = 1632.5
$ws.Cells.Item(113, 10).Value = 1529.875
$ws.Cells.Item(113, 11).Value = 1632.5
$ws.Cells.Item(113, 12).Value = 1529.875
$ws.Cells.Item(113, 13).Value = 537.5
$ws.Cells.Item(113, 14).Value = -5869.875

$ws.Cells.Item(118, 8).Value = 38302
$ws.Cells.Item(118, 10).Value = 38302
$ws.Cells.Item(118, 12).Value = 38302
$ws.Cells.Item(118, 14).Value = -41616


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2240.6904
$ws.Cells.Item(7, 9).Value = 1937.5927
$ws.Cells.Item(7, 10).Value = 2786.2666
$ws.Cells.Item(7, 11).Value = 1937.5927
$ws.Cells.Item(7, 12).Value = 2786.2666
$ws.Cells.Item(7, 13).Value = -1825.5927
$ws.Cells.Item(7, 14).Value = -3010.2666

$ws.Cells.Item(43, 8).Value = 10000
$ws.Cells.Item(43, 10).Value = 10000
$ws.Cells.Item(43, 12).Value = 10000
$ws.Cells.Item(43, 14).Value = -10386

$ws.Cells.Item(56, 8).Value = 8582
$ws.Cells.Item(56, 9).Value = 4166.5
$ws.Cells.Item(56, 10).Value = 12997.5
$ws.Cells.Item(56, 11).Value = 4166.5
$ws.Cells.Item(56, 12).Value = 12997.5
$ws.Cells.Item(56, 13).Value = -3475.5
$ws.Cells.Item(56, 14).Value = -14379.5

$ws.Cells.Item(96, 8).Value = 32798
$ws.Cells.Item(96, 10).Value = 32798
$ws.Cells.Item(96, 12).Value = 32798
$ws.Cells.Item(96, 14).Value = -38290

$ws.Cells.Item(126, 8).Value = 2240.6904
$ws.Cells.Item(126, 9).Value = 1937.5927
$ws.Cells.Item(126, 10).Value = 2786.2666
$ws.Cells.Item(126, 11).Value = 5812.7781
$ws.Cells.Item(126, 12).Value = 8358.799800000001
$ws.Cells.Item(126, 13).Value = -3342.7781
$ws.Cells.Item(126, 14).Value = -13298.7998

$ws.Cells.Item(132, 8).Value = 2339.8813
$ws.Cells.Item(132, 9).Value = 1361.85
$ws.Cells.Item(132, 10).Value = 4398.8945
$ws.Cells.Item(132, 11).Value = 4085.55
$ws.Cells.Item(132, 12).Value = 13196.6835
$ws.Cells.Item(132, 13).Value = -1555.55
$ws.Cells.Item(132, 14).Value = -18256.6835

$ws.Cells.Item(136, 8).Value = 2202.16
$ws.Cells.Item(136, 9).Value = 1412.125
$ws.Cells.Item(136, 10).Value = 2573.9412
$ws.Cells.Item(136, 11).Value = 4236.375
$ws.Cells.Item(136, 12).Value = 7721.823600000001
$ws.Cells.Item(136, 13).Value = -1686.375
$ws.Cells.Item(136, 14).Value = -12821.8236


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(99, 8).Value = 40983
$ws.Cells.Item(99, 10).Value = 41000
$ws.Cells.Item(99, 12).Value = 41000
$ws.Cells.Item(99, 14).Value = -46990

$ws.Cells.Item(105, 8).Value = 37996
$ws.Cells.Item(105, 10).Value = 37996
$ws.Cells.Item(105, 12).Value = 37996
$ws.Cells.Item(105, 14).Value = -44984

$ws.Cells.Item(132, 8).Value = 2963.8572
$ws.Cells.Item(132, 9).Value = 2744.6875
$ws.Cells.Item(132, 10).Value = 3376.4119
$ws.Cells.Item(132, 11).Value = 8234.0625
$ws.Cells.Item(132, 12).Value = 10129.2357
$ws.Cells.Item(132, 13).Value = -5704.0625
$ws.Cells.Item(132, 14).Value = -15189.2357

$ws.Cells.Item(136, 8).Value = 22312.148
$ws.Cells.Item(136, 9).Value = 32931.16
$ws.Cells.Item(136, 10).Value = 1737.8125
$ws.Cells.Item(136, 11).Value = 98793.48000000001
$ws.Cells.Item(136, 12).Value = 5213.4375
$ws.Cells.Item(136, 13).Value = -96243.48000000001
$ws.Cells.Item(136, 14).Value = -10313.4375

